$d = $word.ActiveDocument
$wordNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- Paragraph 1: "On Pilgrimage - September 1977" (Heading1 + single run)
#     -> Title style, split into one run per word/space token
$titleXml = "<w:p $wordNs>" +
              "<w:pPr><w:pStyle w:val='Title'/></w:pPr>" +
              "<w:r><w:t xml:space='preserve'>On</w:t></w:r>" +
              "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
              "<w:r><w:t xml:space='preserve'>Pilgrimage</w:t></w:r>" +
              "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
              "<w:r><w:t xml:space='preserve'>-</w:t></w:r>" +
              "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
              "<w:r><w:t xml:space='preserve'>September</w:t></w:r>" +
              "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
              "<w:r><w:t xml:space='preserve'>1977</w:t></w:r>" +
            "</w:p>"

$p1 = $d.Paragraphs(1)
$p1.Range.InsertXML($titleXml)

# --- Paragraph 2: "By Dorothy Day" (bold run)
#     -> Authors style, "Dorothy Day" split into tokens, "By " dropped
$authorsXml = "<w:p $wordNs>" +
                "<w:pPr><w:pStyle w:val='Authors'/></w:pPr>" +
                "<w:r><w:t xml:space='preserve'>Dorothy</w:t></w:r>" +
                "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
                "<w:r><w:t xml:space='preserve'>Day</w:t></w:r>" +
              "</w:p>"

$p2 = $d.Paragraphs(2)
$p2.Range.InsertXML($authorsXml)
